$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep these columns as literal text so values like "598.80" or "69.449.07"
# are not silently reinterpreted as numbers (which would drop the exact
# formatting / thousands-separated "price" strings used in this sheet).
$ws.Range("D2:E51").NumberFormat = "@"

# --- Row-by-row Price (D) / Volume(1h) (E) refresh ---
$ws.Range("D2").Value = "69.449.07"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "3.499.20"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "598.80"
$ws.Range("E5").Value = "  -2.85%  "
$ws.Range("D6").Value = "193.99"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("D7").Value = "0.620"
$ws.Range("E7").Value = "  -1.67%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "0.200"
$ws.Range("E9").Value = "  -5.73%  "
$ws.Range("D10").Value = "0.643"
$ws.Range("E10").Value = "  -3.15%  "
$ws.Range("D11").Value = "52.92"
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("D12").Value = "0.0000298"
$ws.Range("E12").Value = "  -3.68%  "
$ws.Range("D13").Value = "9.41"
$ws.Range("E13").Value = "  -1.88%  "
$ws.Range("D14").Value = "4.052.39"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").Value = "596.05"
$ws.Range("E15").Value = "  -3.88%  "
$ws.Range("D16").Value = "69.581.89"
$ws.Range("D17").Value = "18.89"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "12.59"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").Value = "0.123"
$ws.Range("E19").Value = "  +2.05%  "
$ws.Range("D20").Value = "3.487.95"
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").Value = "0.978"
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("D22").Value = "17.76"
$ws.Range("E22").Value = "  +3.25%  "
$ws.Range("D23").Value = "5.30"
$ws.Range("E23").Value = "  +4.91%  "
$ws.Range("D24").Value = "102.24"
$ws.Range("E24").Value = "  -6.14%  "
$ws.Range("E25").Value = "  -3.00%  "
$ws.Range("E26").Value = "  -1.19%  "
$ws.Range("D27").Value = "10.75"
$ws.Range("E27").Value = "  -2.15%  "
$ws.Range("D28").Value = "9.44"
$ws.Range("E28").Value = "  -2.78%  "
$ws.Range("D29").Value = "32.88"
$ws.Range("E29").Value = "  -4.29%  "
$ws.Range("D30").Value = "4.27"
$ws.Range("E30").Value = "  +8.75%  "
$ws.Range("D31").Value = "6.94"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").Value = "12.24"
$ws.Range("E32").Value = "  -2.27%  "
$ws.Range("E33").Value = "  -2.69%  "
$ws.Range("D34").Value = "63.07"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").Value = "3.18"
$ws.Range("E35").Value = "  +2.23%  "
$ws.Range("D36").Value = "3.740.30"
$ws.Range("E36").Value = "  +2.11%  "
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "0.0₃0806"
$ws.Range("E38").Value = "  +3.26%  "
$ws.Range("D39").Value = "3.65"
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("D40").Value = "0.387"
$ws.Range("E40").Value = "  -1.56%  "
$ws.Range("D41").Value = "493.03"
$ws.Range("E41").Value = "  -5.28%  "
$ws.Range("D42").Value = "35.85"
$ws.Range("E42").Value = "  -2.33%  "
$ws.Range("E43").Value = "  -4.27%  "
$ws.Range("D44").Value = "0.0447"
$ws.Range("E44").Value = "  -4.97%  "
$ws.Range("D45").Value = "0.139"
$ws.Range("E45").Value = "  -2.69%  "
$ws.Range("E46").Value = "  -4.84%  "
$ws.Range("D47").Value = "3.23"
$ws.Range("E47").Value = "  -2.83%  "
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").Value = "8.37"
$ws.Range("E49").Value = "  -4.68%  "

# --- Rows 50/51 swap places: FLOKI moves above OceanProtocol, both refreshed ---
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").Value = "0.000243"
$ws.Range("E50").Value = "  +1.02%  "

$ws.Range("B51").Value = "OceanProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
$ws.Range("D51").Value = "1.34"
$ws.Range("E51").Value = "  -0.07%  "
